$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns keep their original Text format so that
# values such as "0.999" or "1.00" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '42.787.94'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '2.528.06'
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '310.94'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = '101.02'
$ws.Range("E6").Value = '  +2.63%  '
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("D10").Value = '35.81'
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("D11").Value = '0.0805'
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").Value = '7.34'
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("E13").Value = '  +0.97%  '
$ws.Range("D14").Value = '2.920.18'
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '15.35'
$ws.Range("E15").Value = '  -2.98%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.542.67'
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").Value = '0.814'
$ws.Range("E17").Value = '  -3.04%  '
$ws.Range("D18").Value = '42.765.21'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = '6.67'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '12.41'
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("D22").Value = '69.87'
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").Value = '243.35'
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("D25").Value = '2.03'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '25.47'
$ws.Range("E27").Value = '  -5.55%  '
$ws.Range("E28").Value = '  -2.26%  '
$ws.Range("D29").Value = '10.18'
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").Value = '38.65'
$ws.Range("E30").Value = '  -3.30%  '
$ws.Range("D31").Value = '161.89'
$ws.Range("E31").Value = '  +1.91%  '
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("D33").Value = '2.81'
$ws.Range("E33").Value = '  +9.03%  '
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").Value = '0.0789'
$ws.Range("E35").Value = '  -1.09%  '
$ws.Range("D36").Value = '18.34'
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("E37").Value = '  -6.99%  '
$ws.Range("D38").Value = '3.08'
$ws.Range("E38").Value = '  -6.02%  '
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("D40").Value = '0.117'
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("D42").Value = '22.18'
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("E44").Value = '  +3.70%  '
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("D46").Value = '1.988.42'
$ws.Range("E46").Value = '  -0.26%  '
$ws.Range("D47").Value = '9.10'
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("D48").Value = '2.774.11'
$ws.Range("E48").Value = '  -1.37%  '
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("D50").Value = '79.51'
$ws.Range("E50").Value = '  -2.19%  '
$ws.Range("E51").Value = '  -2.02%  '
